$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = -4
$ws.Range("F6").Value = -5
$ws.Range("F11").Value = -3
$ws.Range("F15").Value = 6
$ws.Range("F18").Value = -7
$ws.Range("F26").Value = 14
$ws.Range("F27").Value = -3
$ws.Range("F29").Value = -3
$ws.Range("F30").Value = -4
